$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows (2-6) to their new values ---

# Row 2 (existing) - update C2 (place) and D2 (timestamp); B2 stays "admin"
$ws.Range("C2").Value = "Cinema hall, floor 1, upper middle section, row D, seat 2"
$ws.Range("D2").Value = 45508.49010777778

# Row 3 (existing) - update C3 (place); B3 stays "admin"; D3 unchanged
$ws.Range("C3").Value = "Cinema hall, floor 1, upper middle section, row H, seat 3"

# Row 4 (existing) - B4 changes from "bogdan.yakupov@nu.edu.kz" to "admin"; C4 and D4 updated
$ws.Range("B4").Value = "admin"
$ws.Range("C4").Value = "Cinema hall, floor 1, upper left section, row H, seat 3"
$ws.Range("D4").Value = 45517.61121153936

# Row 5 (existing) - update C5 (place) and D5 (timestamp); B5 stays "admin"
$ws.Range("C5").Value = "Cinema hall, floor 1, upper left section, row F, seat 2"
$ws.Range("D5").Value = 45509.57961226852

# Row 6 (existing) - B6 changes from "admin" to "bogdan@nu.edu.kz"; C6 and D6 updated
$ws.Range("B6").Value = "bogdan@nu.edu.kz"
$ws.Range("C6").Value = "Cinema hall, floor 1, upper middle section, row F, seat 5"
$ws.Range("D6").Value = 45509.62525357998

# --- Append new rows (7-9), copying formatting from the existing data rows ---

# Row 7 (new)
$ws.Range("A2").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "bogdan.yakupov@nu.edu.kz"
$ws.Range("C7").Value = "Cinema hall, floor 1, upper left section, row F, seat 4"
$ws.Range("D2").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("D7").Value = 45511.72776239251

# Row 8 (new)
$ws.Range("A2").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "admin"
$ws.Range("C8").Value = "Cinema hall, floor 1, upper left section, row E, seat 4"
$ws.Range("D2").Copy()
$ws.Range("D8").PasteSpecial(-4122)
$ws.Range("D8").Value = 45517.61121153936

# Row 9 (new)
$ws.Range("A2").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "admin"
$ws.Range("C9").Value = "Cinema hall, floor 1, upper left section, row A, seat 1"
$ws.Range("D2").Copy()
$ws.Range("D9").PasteSpecial(-4122)
$ws.Range("D9").Value = 45517.61521243056
